# Update crypto price/volume data on Sheet1 to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.871.96"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.603.00"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.97"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.69"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.81"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.063.96"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.830.57"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.86"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.615.03"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.45"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.26"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.09"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.57"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.03"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.906"
$ws.Range("E35").Value = "  +7.74%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.875"
$ws.Range("E37").Value = "  +4.37%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.60"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "281.59"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.599"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0538"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.61"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.946.71"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.02"
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.50"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.93"
$ws.Range("E51").Value = "  -2.73%  "
